# Commit: "Sun, Jun 28, 2020  8:04:39 PM"
#
# Target edit (per the OOXML diff):
#   ppt/slides/slide5.xml - the <a:tbl> inside the graphicFrame on slide 5
#   had its <a:tableStyleId> changed from the custom/local style
#   {30948872-E80D-4221-9D32-9DD1A08AD4A1} (defined in ppt/tableStyles.xml)
#   to the built-in table style {305F7E56-4A5B-4BF6-967A-BD7D1DA53936}
#   (i.e. the table was re-styled from the Table Design gallery).
#
# We look the table up defensively (scan every slide / shape) rather than
# hard-coding slide/shape indices, and re-apply the style through the real
# Table object so PowerPoint (re)writes <a:tableStyleId> exactly like the
# Table Design gallery would.

$OLD_STYLE_ID = "{30948872-E80D-4221-9D32-9DD1A08AD4A1}"
$NEW_STYLE_ID = "{305F7E56-4A5B-4BF6-967A-BD7D1DA53936}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            # Only re-style the table(s) that still carry the old,
            # document-local style id, mirroring the diff precisely
            # (and leaving any already-restyled / unrelated table alone).
            if ($tbl.Style -eq $OLD_STYLE_ID) {
                $tbl.ApplyStyle($NEW_STYLE_ID)
            }
        }
    }
}
